# CIERRE 27 NOV 2023
# - Update the "amount in words" text on the VALES DE INSENTIVOS voucher
#   (08/100 -> 00/100)
# - Move the active selection on that sheet to F8
# - Make the ARQUITECTO sheet the active/selected tab (was VALES DE INSENTIVOS)

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")
$wsArq   = $wb.Worksheets.Item("ARQUITECTO        ")

# Correct the amount-in-words cell (merged A2:D2)
$wsVales.Range("A2").Value = "SEIS  MIL   PESOS 00/100 M.N."

# Leave the cursor on F8 for this sheet before switching away from it
$wsVales.Range("F8").Select()

# Switch the active tab back to ARQUITECTO
$wsArq.Activate()
